$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.345.76'
$ws.Range('E2').Value = '  -2.89%  '
$ws.Range('D3').Value = '1.938.24'
$ws.Range('E3').Value = '  -3.02%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''250.66'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').Value = '''0.7220'
$ws.Range('E6').Value = '  -6.70%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''0.3328'
$ws.Range('E8').Value = '  -4.62%  '
$ws.Range('D9').Value = '''28.27'
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('D10').Value = '''0.07239'
$ws.Range('D11').Value = '''0.8119'
$ws.Range('E11').Value = '  -3.62%  '
$ws.Range('D12').Value = '''0.08110'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').Value = '1.938.93'
$ws.Range('E13').Value = '  -2.94%  '
$ws.Range('D14').Value = '''5.478'
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('D15').Value = '''94.45'
$ws.Range('E15').Value = '  -6.43%  '
$ws.Range('D16').Value = '''15.04'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = '30.361.96'
$ws.Range('E17').Value = '  -2.81%  '
$ws.Range('D18').Value = '''0.000008250'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '''249.19'
$ws.Range('E19').Value = '  -8.78%  '
$ws.Range('D20').Value = '''5.902'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').Value = '2.194.13'
$ws.Range('E21').Value = '  -2.67%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '''1.003'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = '''6.978'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('D25').Value = '''9.753'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('D26').Value = '''163.24'
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').Value = '''2.388'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').Value = '''19.29'
$ws.Range('E28').Value = '  -3.25%  '
$ws.Range('D29').Value = '''0.1320'
$ws.Range('E29').Value = '  -7.68%  '
$ws.Range('D30').Value = '''1.567'
$ws.Range('E30').Value = '  -1.79%  '
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').Value = '''4.439'
$ws.Range('E32').Value = '  -3.64%  '
$ws.Range('D33').Value = '''4.184'
$ws.Range('E33').Value = '  -5.86%  '
$ws.Range('D34').Value = '''0.05210'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').Value = '''1.294'
$ws.Range('D36').Value = '''0.7496'
$ws.Range('E36').Value = '  -4.86%  '
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('D38').Value = '''0.01979'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').Value = '''2.834'
$ws.Range('E39').Value = '  -2.95%  '
$ws.Range('D40').Value = '''81.06'
$ws.Range('E40').Value = '  -0.85%  '
$ws.Range('D41').Value = '''6.451'
$ws.Range('E41').Value = '  -4.25%  '
$ws.Range('D42').Value = '''0.4547'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('D43').Value = '''2.042'
$ws.Range('E43').Value = '  -4.55%  '
$ws.Range('D44').Value = '''0.8483'
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '''101.97'
$ws.Range('E46').Value = '  -2.87%  '
$ws.Range('D47').Value = '''9.794'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').Value = '''7.454'
$ws.Range('E48').Value = '  -3.38%  '
$ws.Range('D49').Value = '''36.82'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '''0.4191'
$ws.Range('E50').Value = '  -3.10%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.06043'
$ws.Range('E51').Value = '  +0.67%  '
